$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 80.59090999999999
$ws.Range("I6").Value = 80.59090999999999
$ws.Range("K6").Value = 241.77273
$ws.Range("M6").Value = -129.77273

$ws.Range("H33").Value = 71786.57000000001
$ws.Range("I33").Value = 83666.836
$ws.Range("K33").Value = 83666.836
$ws.Range("M33").Value = -83437.836

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H99").Value = 2288.6155
$ws.Range("I99").Value = 1839.2
$ws.Range("J99").Value = 3786.6667
$ws.Range("K99").Value = 5517.6
$ws.Range("L99").Value = 11360.0001
$ws.Range("M99").Value = -4019.6
$ws.Range("N99").Value = -14356.0001

$ws.Range("H106").Value = 4915.5293
$ws.Range("I106").Value = 5721
$ws.Range("K106").Value = 5721
$ws.Range("M106").Value = -5090

$ws.Range("H132").Value = 53103.45
$ws.Range("I132").Value = 55694.156
$ws.Range("K132").Value = 167082.468
$ws.Range("M132").Value = -164552.468

$ws.Range("H137").Value = 1391.5278
$ws.Range("I137").Value = 1305.5652
$ws.Range("K137").Value = 3916.6956
$ws.Range("M137").Value = -1366.6956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 137.7
$ws.Range("J5").Value = 122
$ws.Range("L5").Value = 122
$ws.Range("N5").Value = -346

$ws.Range("H45").Value = 2111
$ws.Range("I45").Value = 2000
$ws.Range("J45").Value = 2222
$ws.Range("K45").Value = 2000
$ws.Range("L45").Value = 2222
$ws.Range("M45").Value = -1623
$ws.Range("N45").Value = -2976

$ws.Range("H74").Value = 3462.1292
$ws.Range("J74").Value = 2246.25
$ws.Range("L74").Value = 2246.25
$ws.Range("N74").Value = -3994.25

$ws.Range("H77").Value = 3462.1292
$ws.Range("J77").Value = 2246.25
$ws.Range("L77").Value = 11231.25
$ws.Range("N77").Value = -19967.25

$ws.Range("H102").Value = 2635.8
$ws.Range("I102").Value = 2067.5454
$ws.Range("K102").Value = 2067.5454
$ws.Range("M102").Value = -445.5454

$ws.Range("H122").Value = 4658.75
$ws.Range("I122").Value = 4996.8
$ws.Range("K122").Value = 14990.4
$ws.Range("M122").Value = -12540.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 137.7
$ws.Range("J4").Value = 122
$ws.Range("L4").Value = 122
$ws.Range("N4").Value = -352

$ws.Range("H20").Value = 3008.1292
$ws.Range("I20").Value = 2368.077
$ws.Range("J20").Value = 3470.389
$ws.Range("K20").Value = 2368.077
$ws.Range("L20").Value = 3470.389
$ws.Range("M20").Value = -2121.077
$ws.Range("N20").Value = -3964.389

$ws.Range("H105").Value = 1993.1765
$ws.Range("I105").Value = 1698
$ws.Range("J105").Value = 2701.6
$ws.Range("K105").Value = 1698
$ws.Range("L105").Value = 2701.6
$ws.Range("M105").Value = 49
$ws.Range("N105").Value = -6195.6

$ws.Range("H134").Value = 1722.5454
$ws.Range("I134").Value = 1618.625
$ws.Range("K134").Value = 4855.875
$ws.Range("M134").Value = -2320.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 2948762.5
$ws.Range("J4").Value = 3579140.5
$ws.Range("L4").Value = 3579140.5
$ws.Range("N4").Value = -3579364.5

$ws.Range("H31").Value = 2587.8696
$ws.Range("J31").Value = 2741.5
$ws.Range("L31").Value = 2741.5
$ws.Range("N31").Value = -3331.5

$ws.Range("H34").Value = 2587.8696
$ws.Range("J34").Value = 2741.5
$ws.Range("L34").Value = 2741.5
$ws.Range("N34").Value = -3145.5

$ws.Range("H60").Value = 37000
$ws.Range("I60").Value = 35000
$ws.Range("K60").Value = 35000
$ws.Range("M60").Value = -34489

$ws.Range("H93").Value = 15832.333
$ws.Range("I93").Value = 8749.75
$ws.Range("K93").Value = 8749.75
$ws.Range("M93").Value = -6877.75

$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -39492

$ws.Range("H105").Value = 2023.7037
$ws.Range("I105").Value = 1977.6
$ws.Range("J105").Value = 2600
$ws.Range("K105").Value = 1977.6
$ws.Range("L105").Value = 2600
$ws.Range("M105").Value = -230.5999999999999
$ws.Range("N105").Value = -6094

$ws.Range("H122").Value = 20612.467
$ws.Range("I122").Value = 3683
$ws.Range("J122").Value = 31898.777
$ws.Range("K122").Value = 11049
$ws.Range("L122").Value = 95696.33099999999
$ws.Range("M122").Value = -8599
$ws.Range("N122").Value = -100596.331

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 345434.28
$ws.Range("J4").Value = 234154.16
$ws.Range("L4").Value = 702462.48
$ws.Range("N4").Value = -702686.48

$ws.Range("H46").Value = 20002918
$ws.Range("J46").Value = 4729.3335
$ws.Range("L46").Value = 14188.0005
$ws.Range("N46").Value = -14370.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 17745
$ws.Range("J5").Value = 17745
$ws.Range("L5").Value = 17745
$ws.Range("N5").Value = -17969

$ws.Range("H70").Value = 9003.444
$ws.Range("I70").Value = 9354.666999999999
$ws.Range("J70").Value = 8652.223
$ws.Range("K70").Value = 9354.666999999999
$ws.Range("L70").Value = 8652.223
$ws.Range("M70").Value = -9084.666999999999
$ws.Range("N70").Value = -9192.223

$ws.Range("H73").Value = 9003.444
$ws.Range("I73").Value = 9354.666999999999
$ws.Range("J73").Value = 8652.223
$ws.Range("K73").Value = 9354.666999999999
$ws.Range("L73").Value = 8652.223
$ws.Range("M73").Value = -8418.666999999999
$ws.Range("N73").Value = -10524.223

$ws.Range("H97").Value = 1243.5217
$ws.Range("I97").Value = 1000.1579
$ws.Range("K97").Value = 1000.1579
$ws.Range("M97").Value = -504.1579

$ws.Range("H122").Value = 4838.8
$ws.Range("I122").Value = 6131
$ws.Range("K122").Value = 18393
$ws.Range("M122").Value = -15943

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6935.5
$ws.Range("I2").Value = 3000
$ws.Range("J2").Value = 18742
$ws.Range("K2").Value = 3000
$ws.Range("L2").Value = 18742
$ws.Range("M2").Value = -2888
$ws.Range("N2").Value = -18966

$ws.Range("H68").Value = 5479.9375
$ws.Range("I68").Value = 2590.9092
$ws.Range("J68").Value = 11835.8
$ws.Range("K68").Value = 2590.9092
$ws.Range("L68").Value = 11835.8
$ws.Range("M68").Value = -1841.9092
$ws.Range("N68").Value = -13333.8

$ws.Range("H71").Value = 5479.9375
$ws.Range("I71").Value = 2590.9092
$ws.Range("J71").Value = 11835.8
$ws.Range("K71").Value = 12954.546
$ws.Range("L71").Value = 59179
$ws.Range("M71").Value = -9210.546
$ws.Range("N71").Value = -66667

$ws.Range("H136").Value = 5311.8237
$ws.Range("I136").Value = 5051.0835
$ws.Range("K136").Value = 15153.2505
$ws.Range("M136").Value = -12603.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 41453.727
$ws.Range("J2").Value = 43199.4
$ws.Range("L2").Value = 43199.4
$ws.Range("N2").Value = -43423.4

$ws.Range("H122").Value = 2049.25

$ws.Range("H126").Value = 4154.4
$ws.Range("I126").Value = 3890.6667
$ws.Range("J126").Value = 4550
$ws.Range("K126").Value = 11672.0001
$ws.Range("L126").Value = 13650
$ws.Range("M126").Value = -9202.000100000001
$ws.Range("N126").Value = -18590

$ws.Range("H136").Value = 2131.0667
$ws.Range("I136").Value = 2197.36
$ws.Range("K136").Value = 6592.08
$ws.Range("M136").Value = -4042.08
